$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells are plain text in the source workbook (t="inlineStr").
# Several of the new values look numeric (e.g. "226.60", "32.90", "1.00") and
# Excel would silently convert them to real numbers (dropping significant
# trailing zeros) unless the cell is explicitly formatted as Text first. Only
# apply this to the values that Excel could actually parse as a number; the
# others (e.g. "34.636.44", which has two dots) are never auto-converted, so
# leave their formatting untouched.
$dCellsNeedingTextFormat = @("D5","D8","D10","D11","D13","D18","D19","D21","D23","D25","D26","D27","D29","D30","D32","D33","D36","D37","D40","D44","D45","D46","D49")
foreach ($cellRef in $dCellsNeedingTextFormat) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.636.44"
$ws.Range("D3").Value = "1.793.21"
$ws.Range("D5").Value = "226.60"
$ws.Range("D8").Value = "32.90"
$ws.Range("D10").Value = "0.0693"
$ws.Range("D11").Value = "0.0952"
$ws.Range("D13").Value = "11.09"
$ws.Range("D14").Value = "1.781.85"
$ws.Range("D16").Value = "34.546.96"
$ws.Range("D18").Value = "68.78"
$ws.Range("D19").Value = "248.04"
$ws.Range("D20").Value = "0.0₃0799"
$ws.Range("D21").Value = "11.25"
$ws.Range("D23").Value = "4.17"
$ws.Range("D25").Value = "165.24"
$ws.Range("D26").Value = "7.28"
$ws.Range("D27").Value = "16.56"
$ws.Range("D29").Value = "1.00"
$ws.Range("D30").Value = "4.17"
$ws.Range("D32").Value = "0.0523"
$ws.Range("D33").Value = "1.23"
$ws.Range("D35").Value = "1.427.36"
$ws.Range("D36").Value = "2.58"
$ws.Range("D37").Value = "0.672"
$ws.Range("D40").Value = "85.15"
$ws.Range("D44").Value = "13.61"
$ws.Range("D45").Value = "0.0527"
$ws.Range("D46").Value = "6.13"
$ws.Range("D49").Value = "106.07"

# E-column "Volume(1h)" percentage strings (already text; keep leading/trailing
# double spaces exactly as in the source).
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E8").Value = "  +3.32%  "
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("E21").Value = "  +2.47%  "
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +14.41%  "
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("E36").Value = "  +6.19%  "
$ws.Range("E37").Value = "  +2.48%  "
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +6.22%  "
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("E45").Value = "  +3.76%  "
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -4.66%  "
